$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Eintrittsdatum (row 7) now has a value ---
$ws.Range("B7").Value = "01.01.2024"

# --- widen column A to fit the longer new labels ---
$ws.Columns.Item(1).ColumnWidth = 46.21875

# ===================================================================
# New rows 22-45: "Daten" / "Wert" style two-column block, continuing
# the alternating-fill convention already used by rows 1-21.
# ===================================================================

# Row 22 - Geschlecht (new accent color #1)
$ws.Range("A22:B22").Interior.ThemeColor = 9
$ws.Range("A22").Value = "Geschlecht"

# Row 23 - Mitarbeitertyp (new accent color #2)
$ws.Range("A23:B23").Interior.ThemeColor = 8
$ws.Range("A23").Value = "Mitarbeitertyp"

# Row 24 - Gesellschaft (new accent color #3)
$ws.Range("A24:B24").Interior.ThemeColor = 7
$ws.Range("A24").Value = "Gesellschaft"

# Rows 25-26 - Jobtitel / Erfahrungsstufe (reuse existing style from A2)
$ws.Range("A2").Copy()
$ws.Range("A25:B26").PasteSpecial(-4122)
$ws.Range("A25").Value = "Jobtitel"
$ws.Range("A26").Value = "Erfahrungsstufe"

# Row 27 - Abteilung (reuse existing style from A16)
$ws.Range("A16").Copy()
$ws.Range("A27:B27").PasteSpecial(-4122)
$ws.Range("A27").Value = "Abteilung"

# Row 28 - Wochenarbeitszeit (accent color #1 again)
$ws.Range("A28:B28").Interior.ThemeColor = 9
$ws.Range("A28").Value = "Wochenarbeitszeit"

# Row 29 - Steuerklasse (accent color #2 again)
$ws.Range("A29:B29").Interior.ThemeColor = 8
$ws.Range("A29").Value = "Steuerklasse"

# Rows 30-31 - Tarifbeschaeftigt? / Tarif (accent color #3 again)
$ws.Range("A30:B31").Interior.ThemeColor = 7
$ws.Range("A30").Value = "Tarifbeschaeftigt?"
$ws.Range("A31").Value = "Tarif"

# Rows 32-35 - aussertariflich beschaeftigt? / AT-Grundgehalt / AT-Weihnachtsgeld / AT-Urlaubsgeld
$ws.Range("A2").Copy()
$ws.Range("A32:B35").PasteSpecial(-4122)
$ws.Range("A32").Value = "aussertariflich beschaeftigt?"
$ws.Range("A33").Value = "AT-Grundgehalt"
$ws.Range("A34").Value = "AT-Weihnachtsgeld"
$ws.Range("A35").Value = "AT-Urlaubsgeld"

# Rows 36-41 - privat versichert? / Zuschuss ... / unfallversichert? / rentenversichert?
$ws.Range("A16").Copy()
$ws.Range("A36:B41").PasteSpecial(-4122)
$ws.Range("A36").Value = "privat versichert?"
$ws.Range("A37").Value = "Zuschuss private Krankenversicherung"
$ws.Range("A38").Value = "Zuschuss privater Zusatzbeitrag"
$ws.Range("A39").Value = "Zuschuss private Pflegeversicheurng"
$ws.Range("A40").Value = "unfallversichert?"
$ws.Range("A41").Value = "rentenversichert?"

# Rows 42-45 - gesetzlich versichert? / Mitglied... / wohnhaft Sachsen? (accent color #1 again)
$ws.Range("A42:B45").Interior.ThemeColor = 9
$ws.Range("A42").Value = "gesetzlich versichert?"
$ws.Range("A43").Value = "Mitglied gesetzliche Krankenkasse (Abkürzung)"
$ws.Range("A44").Value = "Mitglied gesetzliche Krankenkasse (vollständiger Name)"
$ws.Range("A45").Value = "wohnhaft Sachsen?"

# Row 46 - trailing empty row, explicit "no fill" formatting only
$ws.Cells.Item(46, 1).Interior.ColorIndex = -4142

# --- move the active selection to reflect where editing left off ---
$ws.Range("E30").Select()
